$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add & remove a hyperlink at D7 to leave residual "Hyperlink" style with no fill
$ws.Hyperlinks.Add($ws.Range("D7"), "http://example.com") | Out-Null
$ws.Range("D7").Hyperlinks.Delete()
$ws.Range("D7").ClearContents()

# 2. Insert a new row at 9
$ws.Rows("9:9").Insert()

# 3. Fix up C column (stays at row 10 instead of shifting to row 11)
$ws.Range("C11").Copy($ws.Range("C10"))
$ws.Range("C11").Clear()

# 4. Set the new row9 content + hyperlink (webroot) BEFORE renaming A8,
#    so that this string gets shared-string index 72
$ws.Hyperlinks.Add($ws.Range("A9"), "http://www.webroot.com/En_US/sites/aff-wsc-29/?ref=cj&rc=2614") | Out-Null
$ws.Range("B9").Value = "IN PROGRESS"

# 5. Edit A8 text (this new string becomes shared-string index 73)
$ws.Range("A8").Value = "Linked-in => add everybody from the Red Poole"

# 6. Update the selection to reflect where the user ended up (A13)
$ws.Range("A13").Select() | Out-Null

Write-Host "done"
